$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 201; existing rows 201-219 shift down to 202-220.
$ws.Rows.Item(201).Insert()

# Populate the newly inserted row 201 with the new weekly record.
$ws.Range("A201").Value = 4
$ws.Range("B201").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C201").Value = "Los Lagos"
$ws.Range("D201").Value = 44617
$ws.Range("D201").NumberFormat = $ws.Range("D202").NumberFormat
$ws.Range("E201").Value = 10
$ws.Range("F201").Value = "Fruta"
$ws.Range("G201").Value = 100104
$ws.Range("H201").Value = "Frutos de pepita"
$ws.Range("I201").Value = 100104005
$ws.Range("J201").Value = "Pera"
$ws.Range("K201").Value = "Packham's Triumph"
$ws.Range("L201").Value = "Primera"
$ws.Range("M201").Value = 500
$ws.Range("N201").Value = 15000
$ws.Range("O201").Value = 16000
$ws.Range("P201").Value = 15500
$ws.Range("Q201").Value = "`$/caja 15 kilos empedrada"
$ws.Range("R201").Value = "Región de O'Higgins"
$ws.Range("S201").Value = 1033
$ws.Range("T201").Value = 15
